# Updated symbol list (coin price/volume refresh) applied via Excel COM interop.
# D (Price) and E (Volume 1h) columns hold numeric-looking text, so NumberFormat
# is forced to "@" (Text) before assignment to avoid Excel auto-converting the
# string into a real number/percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "285.34"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.06%"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.62"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "5.60%"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.913"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.25%"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06491"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.32%"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.216"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.07%"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.341"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "12.10%"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9137"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "4.23%"

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.28%"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06333"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "23.45%"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07638"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.04%"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02976"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.37%"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08963"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.17%"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001607"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.65%"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006536"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.66%"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006025"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.99%"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.461"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.368"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.82%"

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.43%"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3149"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.42%"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1343"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.35%"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.019"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.98%"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1556"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "12.74%"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04472"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.55%"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001188"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.02%"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004329"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "12.16%"

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001180"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-9.18%"

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0001636"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "-15.71%"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04157"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.07%"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006683"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.80%"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1232"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.02%"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002151"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "10.28%"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01178"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.11%"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005372"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.27%"

# Row 46
$ws.Range("B46").Value = "CoinbaseStockToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.01850"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.11%"

# Row 47
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.041"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "20.97%"
